# Adapt column header formatting to respective input file names.
# "_old" headers -> "_FV2304" suffix, "_new" headers -> "_FV2310" suffix,
# then wrap the used range in an Excel Table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (row 1) cell values ---------------------------
$headerMap = @{
    "Segmentname_old"        = "Segmentname_FV2304"
    "Segmentgruppe_old"      = "Segmentgruppe_FV2304"
    "Segment_old"            = "Segment_FV2304"
    "Datenelement_old"       = "Datenelement_FV2304"
    "Segment ID_old"         = "Segment ID_FV2304"
    "Code_old"               = "Code_FV2304"
    "Qualifier_old"          = "Qualifier_FV2304"
    "Beschreibung_old"       = "Beschreibung_FV2304"
    "Bedingungsausdruck_old" = "Bedingungsausdruck_FV2304"
    "Bedingung_old"          = "Bedingung_FV2304"
    "Segmentname_new"        = "Segmentname_FV2310"
    "Segmentgruppe_new"      = "Segmentgruppe_FV2310"
    "Segment_new"            = "Segment_FV2310"
    "Datenelement_new"       = "Datenelement_FV2310"
    "Segment ID_new"         = "Segment ID_FV2310"
    "Code_new"               = "Code_FV2310"
    "Qualifier_new"          = "Qualifier_FV2310"
    "Beschreibung_new"       = "Beschreibung_FV2310"
    "Bedingungsausdruck_new" = "Bedingungsausdruck_FV2310"
    "Bedingung_new"          = "Bedingung_FV2310"
}

$usedRange = $ws.UsedRange
$lastCol = $usedRange.Columns.Count
$lastRow = $usedRange.Rows.Count

for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $current = $cell.Value2
    if ($headerMap.ContainsKey($current)) {
        $cell.Value = $headerMap[$current]
    }
}

# --- 2. Wrap the data range into a native Excel Table --------------------
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 3. Freeze the header row (pane split below row 1) --------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
